$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6").Value = -7.585999999999995
$ws.Range("C7").Value = -13.2583
$ws.Range("E7").Value = 15.9159
$ws.Range("A9").Value = -21.7461
$ws.Range("E10").Value = 16.70810000000001
$ws.Range("C12").Value = -10.8769
$ws.Range("E13").Value = 16.58050000000001
$ws.Range("C14").Value = -14.0309
$ws.Range("D15").Value = -8.802499999999998
$ws.Range("E16").Value = 16.3032
$ws.Range("A18").Value = -22.31080000000002
$ws.Range("A20").Value = -21.19159999999998
$ws.Range("E20").Value = 15.359
$ws.Range("E24").Value = 16.54770000000001
$ws.Range("C26").Value = -12.73060000000001
$ws.Range("A27").Value = -21.55449999999999
$ws.Range("C27").Value = -12.8684
$ws.Range("C29").Value = -11.2258
$ws.Range("D33").Value = -7.504800000000002
$ws.Range("A35").Value = -20.89599999999998
$ws.Range("D35").Value = -8.120399999999997
$ws.Range("C37").Value = -13.99219999999999
$ws.Range("C38").Value = -13.8115
$ws.Range("D38").Value = -9.141099999999989
$ws.Range("E39").Value = 16.3396
$ws.Range("D43").Value = -8.295600000000002
$ws.Range("D44").Value = -7.566000000000007
$ws.Range("D47").Value = -8.309900000000001
$ws.Range("E47").Value = 16.69050000000001
$ws.Range("E48").Value = 17.4746
$ws.Range("C51").Value = -12.5944
$ws.Range("D51").Value = -7.5825
$ws.Range("C52").Value = -11.3136
$ws.Range("E52").Value = 17.38560000000001
$ws.Range("C55").Value = -14.01990000000001
$ws.Range("E56").Value = 16.7347
$ws.Range("D57").Value = -8.456700000000001
$ws.Range("D63").Value = -8.1914
$ws.Range("A69").Value = -21.75759999999999
$ws.Range("C69").Value = -11.2253
$ws.Range("C70").Value = -12.674
$ws.Range("D70").Value = -7.977299999999996
$ws.Range("A76").Value = -19.69829999999999
$ws.Range("A78").Value = -19.93239999999999
$ws.Range("C81").Value = -13.279
$ws.Range("A82").Value = -22.12870000000001
$ws.Range("A83").Value = -21.79350000000001
$ws.Range("C83").Value = -12.80099999999999
$ws.Range("E84").Value = 16.95570000000001
$ws.Range("D88").Value = -7.516799999999995
$ws.Range("A93").Value = -20.9346
$ws.Range("D99").Value = -7.490299999999998
$ws.Range("E100").Value = 16.40730000000001
$ws.Range("E101").Value = 16.81080000000001
$ws.Range("C102").Value = -13.5242
